$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "E85" technology row entirely (E85 is not a ZEV technology and is
# being dropped from this reference policy sheet). Locate it by content so the
# edit is robust even if row positions ever shift.
$e85Cell = $ws.Cells.Find("E85")
$e85Cell.EntireRow.Delete()

# Tighten the ZEV minimum market-share target for 2035-2050 from 100% to 99%.
$ws.Range("T3:W3").Value = 0.99

# Reflect the resulting full-table selection (A1:X8 was the last used range
# before the deleted row collapsed it to A1:X7).
$ws.Range("A1:X8").Select() | Out-Null
